$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.733.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.904.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.56%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5206"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.66%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3778"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.32%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07246"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.98%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.36%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9030"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.51%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07663"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.33%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.904.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.50%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.448"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.29%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.19%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.12%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008710"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.26%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9995"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.17%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "27.764.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.50%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.140"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.49%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.156.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.61%  "

# Row 23
$ws.Range("E23").Value = "  +0.95%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.631"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.29%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.869"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.73%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.30%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.160"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.37%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.52%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.854"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.13%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09026"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.02%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.186"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.45%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.837"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.62%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.232"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.65%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7808"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.37%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02092"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.66%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.597"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.58%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.084"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.29%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.093"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.29%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5567"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.91%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.05286"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.12%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.726"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.07%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "114.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.07%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.523"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.04%  "

# Row 45
$ws.Range("E45").Value = "  -0.02%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4815"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.68%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.39%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9996"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.12%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.615"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.12%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.67%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.86%  "
